# Re-sort the worksheet tabs: "总计" (summary) should come first, followed
# by "2022-Q2" (detail). The data on each sheet is left untouched - this is
# purely a reorder of the sheet tabs.

$wb = $excel.ActiveWorkbook

# Move the summary sheet so it becomes the first tab in the workbook.
$wb.Worksheets.Item("总计").Move($wb.Worksheets.Item(1))

# Keep "2022-Q2" as the active/selected sheet, same as before the reorder.
# (Re-fetch by name rather than reusing a pre-move reference, since the
# sheet's position shifted.)
$wb.Worksheets.Item("2022-Q2").Activate()

Write-Output ("Sheet order: " + (($wb.Worksheets | ForEach-Object { $_.Name }) -join ", "))
